# "working trials, updated sheet"
#
# Update the Trials&Showdowns worksheet:
#  - Column A values lose the trailing " Trial" suffix (e.g. "Dark Trial" -> "Dark")
#  - Column B (monster names) is left untouched
#  - The sheet view selection is updated to cover the whole data column
#    and the window is scrolled down a bit

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dim = $ws.UsedRange
$lastRow = $dim.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = [string]$cell.Value2
    if ($val -match '^(.*) Trial$') {
        $cell.Value2 = $Matches[1]
    }
}

# Update selection / scroll position to match the saved view state
# (select the full data column and scroll the window so row 10 is at the top)
$ws.Range("A1:A31").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
